$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "USURIO : " -> "USUARIO : "
$ws.Range("A5").Value = "USUARIO : "

# Add accent: "CATEGORIA" -> "CATEGORÍA"
$ws.Range("A8").Value = "CATEGORÍA"

# Update the saved selection to A8
$ws.Range("A8").Select()
